$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove leftover hidden chart-related defined names (_xlchart.*)
# ---------------------------------------------------------------------------
while ($wb.Names.Count -gt 0) {
    $wb.Names.Item(1).Delete()
}

# ---------------------------------------------------------------------------
# 2. Update the existing "deve essere" hint strings (precedente / posizione 8)
#    from "- addX" to "+ addX" wording. We first move the affected cells to a
#    placeholder value so that the old shared-string entries become orphaned
#    (and get garbage collected on save), then set the real new text so that
#    the shared string table is rebuilt in the same order as the target file.
# ---------------------------------------------------------------------------
$ws.Range("J6").Value = "__TMP1__"
$ws.Range("J8").Value = "__TMP1__"
$ws.Range("J15").Value = "__TMP1__"
$ws.Range("J10").Value = "__TMP2__"
$ws.Range("J11").Value = "__TMP2__"
$ws.Range("J12").Value = "__TMP2__"
$ws.Range("J13").Value = "__TMP2__"

$ws.Range("J6").Value = "deve essere 9 + addY precedente + addX"
$ws.Range("J8").Value = "deve essere 9 + addY precedente + addX"
$ws.Range("J15").Value = "deve essere 9 + addY precedente + addX"
$ws.Range("J10").Value = "deve essere 9 + addY di posizione 8 + addX o al minimo 9"
$ws.Range("J11").Value = "deve essere 9 + addY di posizione 8 + addX o al minimo 9"
$ws.Range("J12").Value = "deve essere 9 + addY di posizione 8 + addX o al minimo 9"
$ws.Range("J13").Value = "deve essere 9 + addY di posizione 8 + addX o al minimo 9"

# ---------------------------------------------------------------------------
# 3. Fix a couple of formulas in the existing part-A table
# ---------------------------------------------------------------------------
$ws.Range("O2").Formula = "=B13 - D13 - E2"
$ws.Range("O11").Formula = "=MIN(9,9+E4+D11)"
$ws.Range("O12").Formula = "=MIN(9,9+E3+D12)"
$ws.Range("O13").Formula = "=MIN(9,9+E2+D13)"

# ---------------------------------------------------------------------------
# 4. Add the new "part B" table (rows 18-31)
# ---------------------------------------------------------------------------

# --- row 18 ---------------------------------------------------------------
$ws.Range("A18").Value = 1
$ws.Range("B18").Value = 1
$ws.Range("C18").Value = 1
$ws.Range("D18").Value = 14
$ws.Range("E18").Value = 12
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Formula = "=B18+E18"
$ws.Range("O18").Formula = "=B29 - D29 - E18"

# --- row 19 ---------------------------------------------------------------
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = 6
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 8
$ws.Range("F19").Formula = "=MOD(H18,26)"
$ws.Range("G19").Formula = "=F19 + D19"
$ws.Range("O19").Formula = "=B28 - D28 - E19"

# --- row 20 ---------------------------------------------------------------
$ws.Range("A20").Value = 3
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 7
$ws.Range("O20").Formula = "=B27 - D27 - E20"

# --- row 21 ---------------------------------------------------------------
$ws.Range("A21").Value = 4
$ws.Range("B21").Value = 8
$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 14
$ws.Range("E21").Value = 4
$ws.Range("O21").Formula = "=B22-D22-E21"

# H19:H21 is a shared formula "IF(G.=B.,TRUNC(H(prev)/26),H(prev)*26+B.+E.)"
$ws.Range("H19:H21").FormulaR1C1 = "=IF(RC[-1]=RC[-6],TRUNC(R[-1]C/26),R[-1]C*26+RC[-6]+RC[-3])"

# --- row 22 -----------------------------------------------------------------
$ws.Range("A22").Value = 5
$ws.Range("B22").Value = 1
$ws.Range("C22").Value = 26
$ws.Range("D22").Value = -11
$ws.Range("E22").Value = 4
$ws.Range("H22").Formula = "=IF(G22=B22,TRUNC(H21/26),H21 * 26 + B22 + E22)"
$ws.Range("J22").Value = "deve essere il massimo fra 1 e 1 + addY precedente + addX"
$ws.Range("O22").Formula = "=MAX(1,1+E21+D22)"

# --- row 23 ---------------------------------------------------------------
$ws.Range("A23").Value = 6
$ws.Range("B23").Value = 1
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 12
$ws.Range("E23").Value = 1
$ws.Range("O23").Formula = "=B24-D24-E23"

# F20:F31 shared formula "MOD(H(prev),26)"
$ws.Range("F20:F31").FormulaR1C1 = "=MOD(R[-1]C[2],26)"
# G20:G31 shared formula "F.+D."
$ws.Range("G20:G31").FormulaR1C1 = "=RC[-1]+RC[-3]"

# --- row 24 ---------------------------------------------------------------
$ws.Range("A24").Value = 7
$ws.Range("B24").Value = 1
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = -1
$ws.Range("E24").Value = 10
$ws.Range("J24").Value = "deve essere il massimo fra 1 e 1 + addY precedente + addX"
$ws.Range("O24").Formula = "=MAX(1,1+E23+D24)"

# --- row 25 ---------------------------------------------------------------
$ws.Range("A25").Value = 8
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 10
$ws.Range("E25").Value = 8
$ws.Range("O25").Formula = "=B26-D26-E25"

# --- row 26 ---------------------------------------------------------------
$ws.Range("A26").Value = 9
$ws.Range("B26").Value = 6
$ws.Range("C26").Value = 26
$ws.Range("D26").Value = -3
$ws.Range("E26").Value = 12
$ws.Range("J26").Value = "deve essere il massimo fra 1 e 1 + addY di posizione 8 + addX"
$ws.Range("O26").Formula = "=MAX(1,1+E25+D26)"

# --- row 27 ---------------------------------------------------------------
$ws.Range("A27").Value = 10
$ws.Range("B27").Value = 4
$ws.Range("C27").Value = 26
$ws.Range("D27").Value = -4
$ws.Range("E27").Value = 10
$ws.Range("J27").Value = "deve essere il massimo fra 1 e 1 + addY di posizione 8 + addX"
$ws.Range("O27").Formula = "=MAX(1,1+E20+D27)"

# --- row 28 ---------------------------------------------------------------
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = 1
$ws.Range("C28").Value = 26
$ws.Range("D28").Value = -13
$ws.Range("E28").Value = 15
$ws.Range("J28").Value = "deve essere il massimo fra 1 e 1 + addY di posizione 8 + addX"
$ws.Range("O28").Formula = "=MAX(1,1+E19+D28)"

# --- row 29 ---------------------------------------------------------------
$ws.Range("A29").Value = 12
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = 26
$ws.Range("D29").Value = -8
$ws.Range("E29").Value = 4
$ws.Range("J29").Value = "deve essere il massimo fra 1 e 1 + addY di posizione 8 + addX"
$ws.Range("O29").Formula = "=MAX(1,1+E18+D29)"

# --- row 30 ---------------------------------------------------------------
$ws.Range("A30").Value = 13
$ws.Range("B30").Value = 2
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 13
$ws.Range("E30").Value = 10
$ws.Range("O30").Formula = "=B31-D31-E30"

# --- row 31 ---------------------------------------------------------------
$ws.Range("A31").Value = 14
$ws.Range("B31").Value = 1
$ws.Range("C31").Value = 26
$ws.Range("D31").Value = -11
$ws.Range("E31").Value = 9
$ws.Range("J31").Value = "deve essere il massimo fra 1 e 1 + addY precedente + addX"
$ws.Range("O31").Formula = "=MAX(1,1+E30+D31)"

# H23:H31 shared formula "IF(G.=B.,TRUNC(H(prev)/26),H(prev)*26+B.+E.)"
$ws.Range("H23:H31").FormulaR1C1 = "=IF(RC[-1]=RC[-6],TRUNC(R[-1]C/26),R[-1]C*26+RC[-6]+RC[-3])"

# ---------------------------------------------------------------------------
# 5. Cosmetics: column J width, selection
# ---------------------------------------------------------------------------
$ws.Columns("J").ColumnWidth = 10

$ws.Range("O13").Select()
